$wb = $excel.ActiveWorkbook

# -------------------------------------------------------------------
# "Extent of Contamination" sheet - data edits (BOTE input updates)
# -------------------------------------------------------------------
$ws = $wb.Worksheets.Item("Extent of Contamination")

# Row 4: Parameter 1 value updated
$ws.Range("G4").Value = 747.86950000000002

# Row 5: Distribution type changed from Constant to Uniform, with new
# Min/Max parameters. Pull the border/box formatting used by the other
# "last row of a record" cells (e.g. F9/G9/H9, which already use that
# style) so the boxed-cell look carries over correctly.
$ws.Range("F9:H9").Copy()
$ws.Range("F5:H5").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F5").Value = "Uniform"
$ws.Range("G5").Value = 3.0319702616531905
$ws.Range("H5").Value = 8.0319702616531909

# Row 6: Parameter 1 value reset to 0
$ws.Range("G6").Value = 0

# Row 7: Distribution type changed from Uniform back to Constant; second
# parameter cleared. G7 previously had no explicit cell style -- copy the
# "continuation row" box style from H7 (already s=12) so it matches.
$ws.Range("H7").Copy()
$ws.Range("G7").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("F7").Value = "Constant"
$ws.Range("G7").Value = -1
$ws.Range("H7").ClearContents()

# Row 9: Parameter 1 value updated
$ws.Range("G9").Value = 1

# Rows 21-32: redistributed Parameter 1 percentages among indoor surface
# breakout rows
$ws.Range("G21").Value = 0
$ws.Range("G22").Value = 0
$ws.Range("G24").Value = 0
$ws.Range("G25").Value = 0
$ws.Range("G27").Value = 0.5
$ws.Range("G28").Value = 0.125
$ws.Range("G29").Value = 0.0625
$ws.Range("G30").Value = 0.0625
$ws.Range("G31").Value = 0.125
$ws.Range("G32").Value = 0.125

# Conditional formatting range now covers the full contiguous block
# A2:L32 instead of the earlier carved-out ranges (A2:L6 A8:L32 A7:F7
# H7:L7) -- the two existing rules just get re-applied to the new range.
$fcs = $ws.Cells.Item(2, 1).FormatConditions
$fcs.Item(1).ModifyAppliesToRange($ws.Range("A2:L32"))
$fcs.Item(2).ModifyAppliesToRange($ws.Range("A2:L32"))

# Restore the active selection to match the saved view
$ws.Range("G7").Select()
